# Update the Correspond Handoff/Handback Datetime values on the
# zh-cn and de-de report sheets (row 2) to reflect the new report
# generation timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 12:41:09"
$wsZhCn.Range("H2").Value = "2016-03-22 12:42:12"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 12:41:17"
$wsDeDe.Range("H2").Value = "2016-03-22 12:42:27"
